# Apply "Select is changed" edit to the pageObjects sheet.
#
# Summary of the change (derived from the OOXML diff):
#   - Row 6 (eightDigitValidation): elementType changes from "assert" to "list",
#     and a new additionalXpath value is added in column D, equal to the xpath
#     already present in column C for that row.
#   - Row 7 (upperCaseValidation): elementType changes from "assert" to "select".
#   - Row 6's height grows (to accommodate the extra content) from 17.35 to 48.15.
#   - The active selection on the sheet moves from D3 to C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# elementType for eightDigitValidation row: assert -> list
$ws.Range("B6").Value = "list"

# additionalXpath for eightDigitValidation row: new cell, same xpath as column C
$ws.Range("D6").Value = $ws.Range("C6").Value()

# elementType for upperCaseValidation row: assert -> select
$ws.Range("B7").Value = "select"

# Row 6 grows taller to fit the newly added additionalXpath column
$ws.Rows.Item(6).RowHeight = 48.15

# Move the active cell/selection to C8
$ws.Range("C8").Select()
